# Append "phase 4" experiment rows (id_response 20-31) to the Qualtrics
# sign-in sheet, extending the used range from A1:K20 to A1:K32.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple: id_response, start_date, end_date, status, progress,
#             duration_s, finished, recorded_date, q1_3, q1_4, q2_2
$data = @(
    ,(20, 44523.78958333333, 44523.79026620371, 'IP Address', 100, 58,  1, 44523.79026620371, 'Yes. I am 18 years of age or older.', 'Yes. I have read the consent form and this response will serve as my consent to participate in the research study.', '1pogus')
    ,(21, 44523.79092592593, 44523.79208333333, 'IP Address', 100, 99,  1, 44523.79208333333, 'Yes. I am 18 years of age or older.', 'Yes. I have read the consent form and this response will serve as my consent to participate in the research study.', '1dabec')
    ,(22, 44523.79164351852, 44523.7924537037,  'IP Address', 100, 70,  1, 44523.7924537037,  'Yes. I am 18 years of age or older.', 'Yes. I have read the consent form and this response will serve as my consent to participate in the research study.', '4lianj')
    ,(23, 44523.79159722223, 44523.79251157407, 'IP Address', 100, 78,  1, 44523.79251157407, 'Yes. I am 18 years of age or older.', 'Yes. I have read the consent form and this response will serve as my consent to participate in the research study.', '3hostc')
    ,(24, 44523.7925462963,  44523.79327546296, 'IP Address', 100, 63,  1, 44523.79327546296, 'Yes. I am 18 years of age or older.', 'Yes. I have read the consent form and this response will serve as my consent to participate in the research study.', '0LIBOL')
    ,(25, 44523.79273148148, 44523.79327546296, 'IP Address', 100, 47,  1, 44523.79328703704, 'Yes. I am 18 years of age or older.', 'Yes. I have read the consent form and this response will serve as my consent to participate in the research study.', '1sawes')
    ,(26, 44523.79199074074, 44523.79365740741, 'IP Address', 100, 144, 1, 44523.79365740741, 'Yes. I am 18 years of age or older.', 'Yes. I have read the consent form and this response will serve as my consent to participate in the research study.', '2dunic')
    ,(27, 44523.7921412037,  44523.79421296297, 'IP Address', 100, 178, 1, 44523.79421296297, 'Yes. I am 18 years of age or older.', 'Yes. I have read the consent form and this response will serve as my consent to participate in the research study.', '1bimil')
    ,(28, 44523.79256944444, 44523.79449074074, 'IP Address', 100, 165, 1, 44523.79449074074, 'Yes. I am 18 years of age or older.', 'Yes. I have read the consent form and this response will serve as my consent to participate in the research study.', '1lacat')
    ,(29, 44523.79333333334, 44523.79483796297, 'IP Address', 100, 130, 1, 44523.79483796297, 'Yes. I am 18 years of age or older.', 'Yes. I have read the consent form and this response will serve as my consent to participate in the research study.', '3mihar')
    ,(30, 44523.7950462963,  44523.79575231481, 'IP Address', 100, 61,  1, 44523.79575231481, 'Yes. I am 18 years of age or older.', 'Yes. I have read the consent form and this response will serve as my consent to participate in the research study.', '2micar')
    ,(31, 44523.79260416667, 44523.79762731482, 'IP Address', 100, 434, 1, 44523.79762731482, 'Yes. I am 18 years of age or older.', 'Yes. I have read the consent form and this response will serve as my consent to participate in the research study.', '1davec')
)

$startRow = 21
for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $r = $startRow + $idx
    $item = $data[$idx]

    $ws.Cells.Item($r, 1).Value2  = $item[0]   # A: id_response
    $ws.Cells.Item($r, 2).Value2  = $item[1]   # B: start_date
    $ws.Cells.Item($r, 3).Value2  = $item[2]   # C: end_date
    $ws.Cells.Item($r, 4).Value2  = $item[3]   # D: status
    $ws.Cells.Item($r, 5).Value2  = $item[4]   # E: progress
    $ws.Cells.Item($r, 6).Value2  = $item[5]   # F: duration_s
    $ws.Cells.Item($r, 7).Value2  = [bool]$item[6]  # G: finished
    $ws.Cells.Item($r, 8).Value2  = $item[7]   # H: recorded_date
    $ws.Cells.Item($r, 9).Value2  = $item[8]   # I: q1_3
    $ws.Cells.Item($r, 10).Value2 = $item[9]   # J: q1_4
    $ws.Cells.Item($r, 11).Value2 = $item[10]  # K: q2_2
}
